$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on columns D and G so numeric-looking values stay as text
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value2 = "242.28"
$ws.Range("G2").Value2 = "19"

# Row 3
$ws.Range("D3").Value2 = "21.56"
$ws.Range("G3").Value2 = "19"

# Row 4
$ws.Range("G4").Value2 = "19"

# Row 5
$ws.Range("D5").Value2 = "0.05610"
$ws.Range("G5").Value2 = "19"

# Row 6
$ws.Range("D6").Value2 = "3.380"
$ws.Range("G6").Value2 = "19"

# Row 7
$ws.Range("D7").Value2 = "6.381"
$ws.Range("G7").Value2 = "19"

# Row 8
$ws.Range("D8").Value2 = "0.8071"
$ws.Range("G8").Value2 = "19"

# Row 9
$ws.Range("D9").Value2 = "0.8918"
$ws.Range("G9").Value2 = "19"

# Row 10
$ws.Range("B10").Value2 = "WazirX"
$ws.Range("C10").Value2 = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value2 = "0.1424"
$ws.Range("E10").Value2 = "9WazirXWRX"
$ws.Range("G10").Value2 = "19"

# Row 11
$ws.Range("B11").Value2 = "MandalaExchangeToken"
$ws.Range("C11").Value2 = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value2 = "0.07315"
$ws.Range("E11").Value2 = "10MandalaExchangeTokenMDX"
$ws.Range("G11").Value2 = "19"

# Row 12
$ws.Range("B12").Value2 = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value2 = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value2 = "0.03232"
$ws.Range("E12").Value2 = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Value2 = "19"

# Row 13
$ws.Range("B13").Value2 = "BitrueCoin"
$ws.Range("C13").Value2 = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value2 = "0.03049"
$ws.Range("E13").Value2 = "12BitrueCoinBTR"
$ws.Range("G13").Value2 = "19"

# Row 14
$ws.Range("B14").Value2 = "BitMartToken"
$ws.Range("C14").Value2 = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value2 = "0.09270"
$ws.Range("E14").Value2 = "13BitMartTokenBMX"
$ws.Range("G14").Value2 = "19"

# Row 15
$ws.Range("B15").Value2 = "MCDex"
$ws.Range("C15").Value2 = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value2 = "3.578"
$ws.Range("E15").Value2 = "14MCDexMCB"
$ws.Range("G15").Value2 = "19"

# Row 16
$ws.Range("B16").Value2 = "BitForexToken"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value2 = "0.001624"
$ws.Range("E16").Value2 = "15BitForexTokenBF"
$ws.Range("G16").Value2 = "19"

# Row 17
$ws.Range("B17").Value2 = "CoinExToken"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value2 = "0.04698"
$ws.Range("E17").Value2 = "16CoinExTokenCET"
$ws.Range("G17").Value2 = "19"

# Row 18
$ws.Range("B18").Value2 = "One"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value2 = "0.0005812"
$ws.Range("E18").Value2 = "17OneONE"
$ws.Range("G18").Value2 = "19"

# Row 19
$ws.Range("D19").Value2 = "0.006363"
$ws.Range("G19").Value2 = "19"

# Row 20
$ws.Range("D20").Value2 = "0.004977"
$ws.Range("G20").Value2 = "19"

# Row 21
$ws.Range("B21").Value2 = "BitKan"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value2 = "0.001043"
$ws.Range("E21").Value2 = "20BitKanKAN"
$ws.Range("G21").Value2 = "19"

# Row 22
$ws.Range("B22").Value2 = "NitroEx"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value2 = "0.0001502"
$ws.Range("E22").Value2 = "21NitroExNTX"
$ws.Range("G22").Value2 = "19"

# Row 23
$ws.Range("B23").Value2 = "UpBots"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D23").Value2 = "0.0003107"
$ws.Range("E23").Value2 = "22UpBotsUBXT"
$ws.Range("G23").Value2 = "19"

# Row 24
$ws.Range("D24").Value2 = "3.768"
$ws.Range("G24").Value2 = "19"

# Row 25
$ws.Range("D25").Value2 = "2.098"
$ws.Range("G25").Value2 = "19"

# Row 26
$ws.Range("D26").Value2 = "0.3252"
$ws.Range("G26").Value2 = "19"

# Row 27
$ws.Range("G27").Value2 = "19"

# Row 28
$ws.Range("G28").Value2 = "19"

# Row 29
$ws.Range("G29").Value2 = "19"

# Row 30
$ws.Range("G30").Value2 = "19"

# Row 31
$ws.Range("G31").Value2 = "19"

# Row 32
$ws.Range("G32").Value2 = "19"

# Row 33
$ws.Range("G33").Value2 = "19"

# Row 34
$ws.Range("G34").Value2 = "19"

# Row 35
$ws.Range("G35").Value2 = "19"

# Row 36
$ws.Range("G36").Value2 = "19"

# Row 37
$ws.Range("G37").Value2 = "19"

# Row 38
$ws.Range("G38").Value2 = "19"

# Row 39
$ws.Range("G39").Value2 = "19"

# Row 40
$ws.Range("D40").Value2 = "0.03889"
$ws.Range("G40").Value2 = "19"

# Row 41
$ws.Range("D41").Value2 = "0.006970"
$ws.Range("G41").Value2 = "19"

# Row 42
$ws.Range("D42").Value2 = "0.1033"
$ws.Range("G42").Value2 = "19"

# Row 43
$ws.Range("D43").Value2 = "0.002924"
$ws.Range("G43").Value2 = "19"

# Row 44
$ws.Range("D44").Value2 = "0.007551"
$ws.Range("G44").Value2 = "19"

# Row 45
$ws.Range("D45").Value2 = "0.00005952"
$ws.Range("G45").Value2 = "19"

# Row 46
$ws.Range("G46").Value2 = "19"

# Row 47
$ws.Range("D47").Value2 = "0.0005511"
$ws.Range("G47").Value2 = "19"

# Row 48
$ws.Range("D48").Value2 = "0.6839"
$ws.Range("G48").Value2 = "19"

# Row 49
$ws.Range("D49").Value2 = "0.05902"
$ws.Range("E49").Value2 = "48BOLOBOLOBestin24h"
$ws.Range("G49").Value2 = "19"

# Row 50
$ws.Range("D50").Value2 = "0.00002104"
$ws.Range("G50").Value2 = "19"

# Row 51
$ws.Range("G51").Value2 = "19"
